$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'37.310.69"
$ws.Range("E2").Value = "'  +2.94%  "
$ws.Range("D3").Value = "'2.060.38"
$ws.Range("E3").Value = "'  +5.79%  "
$ws.Range("E4").Value = "'  -0.08%  "
$ws.Range("D5").Value = "'235.63"
$ws.Range("E5").Value = "'  +4.18%  "
$ws.Range("D6").Value = "'0.615"
$ws.Range("E6").Value = "'  +4.75%  "
$ws.Range("D7").Value = "'57.69"
$ws.Range("E7").Value = "'  +10.22%  "
$ws.Range("E8").Value = "'  -0.11%  "
$ws.Range("E9").Value = "'  +5.48%  "
$ws.Range("D10").Value = "'57.59"
$ws.Range("E10").Value = "'  +1.87%  "
$ws.Range("D11").Value = "'0.0759"
$ws.Range("E11").Value = "'  +4.66%  "
$ws.Range("E12").Value = "'  +4.98%  "
$ws.Range("D13").Value = "'2.366.79"
$ws.Range("E13").Value = "'  +5.81%  "
$ws.Range("D14").Value = "'14.27"
$ws.Range("E14").Value = "'  +5.11%  "
$ws.Range("D15").Value = "'20.84"
$ws.Range("E15").Value = "'  +8.71%  "
$ws.Range("E16").Value = "'  +5.55%  "
$ws.Range("D17").Value = "'5.18"
$ws.Range("E17").Value = "'  +5.57%  "
$ws.Range("D18").Value = "'2.076.67"
$ws.Range("E18").Value = "'  +6.12%  "
$ws.Range("D19").Value = "'37.505.22"
$ws.Range("E19").Value = "'  +3.58%  "
$ws.Range("D20").Value = "'6.11"
$ws.Range("E20").Value = "'  +24.88%  "
$ws.Range("D21").Value = "'68.38"
$ws.Range("E21").Value = "'  +3.28%  "
$ws.Range("D22").Value = "'0.0₃0808"
$ws.Range("E22").Value = "'  +3.71%  "
$ws.Range("D23").Value = "'224.51"
$ws.Range("E23").Value = "'  +3.39%  "
$ws.Range("E24").Value = "'  -0.11%  "
$ws.Range("D25").Value = "'2.45"
$ws.Range("E25").Value = "'  +6.91%  "
$ws.Range("E26").Value = "'  +3.54%  "
$ws.Range("D27").Value = "'163.62"
$ws.Range("E27").Value = "'  +2.77%  "
$ws.Range("D28").Value = "'8.83"
$ws.Range("E28").Value = "'  +6.54%  "
$ws.Range("E29").Value = "'  +10.62%  "
$ws.Range("E30").Value = "'  +9.33%  "
$ws.Range("D31").Value = "'19.17"
$ws.Range("E31").Value = "'  +4.19%  "
$ws.Range("D32").Value = "'0.119"
$ws.Range("E32").Value = "'  +3.66%  "
$ws.Range("D33").Value = "'2.63"
$ws.Range("E33").Value = "'  +18.85%  "
$ws.Range("D34").Value = "'0.0626"
$ws.Range("E34").Value = "'  +6.02%  "
$ws.Range("E35").Value = "'  +4.78%  "
$ws.Range("D36").Value = "'4.45"
$ws.Range("E36").Value = "'  +7.79%  "
$ws.Range("D37").Value = "'1.79"
$ws.Range("E37").Value = "'  +0.71%  "
$ws.Range("E38").Value = "'  -0.10%  "
$ws.Range("D39").Value = "'3.33"
$ws.Range("E39").Value = "'  +7.79%  "
$ws.Range("E40").Value = "'  +17.89%  "
$ws.Range("D41").Value = "'2.99"
$ws.Range("E41").Value = "'  -0.23%  "
$ws.Range("B42").Value = "'Cronos"
$ws.Range("C42").Value = "'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D42").Value = "'0.0953"
$ws.Range("E42").Value = "'  +11.94%  "
$ws.Range("B43").Value = "'FTXToken"
$ws.Range("C43").Value = "'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D43").Value = "'4.41"
$ws.Range("E43").Value = "'  +28.31%  "
$ws.Range("D44").Value = "'1.466.36"
$ws.Range("E44").Value = "'  +5.61%  "
$ws.Range("D45").Value = "'95.16"
$ws.Range("E45").Value = "'  +12.17%  "
$ws.Range("E46").Value = "'  +7.32%  "
$ws.Range("D47").Value = "'16.05"
$ws.Range("E47").Value = "'  +11.61%  "
$ws.Range("E48").Value = "'  +7.04%  "
$ws.Range("E49").Value = "'  +6.02%  "
$ws.Range("D50").Value = "'7.24"
$ws.Range("E50").Value = "'  +9.47%  "
$ws.Range("E51").Value = "'  +2.70%  "
